$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for numeric-looking Price values so Excel
# does not silently convert them to real numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '30.433.47'
$ws.Range("E2").Value = '  +1.73%  '

$ws.Range("D3").Value = '1.688.06'
$ws.Range("E3").Value = '  +3.66%  '

$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.40%  '

$ws.Range("D5").Value = '221.04'
$ws.Range("E5").Value = '  +3.17%  '

$ws.Range("E6").Value = '  +0.47%  '

$ws.Range("D7").Value = '0.996'
$ws.Range("E7").Value = '  -0.35%  '

$ws.Range("D8").Value = '30.35'
$ws.Range("E8").Value = '  +1.91%  '

$ws.Range("E9").Value = '  +2.22%  '

$ws.Range("D10").Value = '0.0623'
$ws.Range("E10").Value = '  +1.78%  '

$ws.Range("D11").Value = '0.0901'
$ws.Range("E11").Value = '  -1.48%  '

$ws.Range("D12").Value = '1.927.78'
$ws.Range("E12").Value = '  +3.53%  '

$ws.Range("D13").Value = '10.78'
$ws.Range("E13").Value = '  +16.98%  '

$ws.Range("E14").Value = '  +8.59%  '

$ws.Range("D15").Value = '1.673.16'
$ws.Range("E15").Value = '  +2.61%  '

$ws.Range("D16").Value = '4.00'
$ws.Range("E16").Value = '  +3.50%  '

$ws.Range("D17").Value = '30.401.28'
$ws.Range("E17").Value = '  +1.63%  '

$ws.Range("D18").Value = '65.86'
$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("D19").Value = '246.53'
$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("E20").Value = '  +2.36%  '

$ws.Range("D21").Value = '0.997'
$ws.Range("E21").Value = '  -0.32%  '

$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '10.21'
$ws.Range("E22").Value = '  +6.34%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '4.30'
$ws.Range("E23").Value = '  +3.85%  '

$ws.Range("D24").Value = '2.21'
$ws.Range("E24").Value = '  +4.17%  '

$ws.Range("D25").Value = '158.72'
$ws.Range("E25").Value = '  -0.17%  '

$ws.Range("D26").Value = '15.90'
$ws.Range("E26").Value = '  +1.29%  '

$ws.Range("E27").Value = '  +0.56%  '

$ws.Range("E28").Value = '  +2.76%  '

$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("E30").Value = '  +2.58%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '3.51'
$ws.Range("E31").Value = '  +4.37%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  +1.25%  '

$ws.Range("D33").Value = '3.31'
$ws.Range("E33").Value = '  +3.56%  '

$ws.Range("D34").Value = '1.511.24'
$ws.Range("E34").Value = '  +5.82%  '

$ws.Range("E35").Value = '  +5.55%  '

$ws.Range("E36").Value = '  -0.04%  '

$ws.Range("E37").Value = '  +5.12%  '

$ws.Range("D38").Value = '2.75'
$ws.Range("E38").Value = '  -4.02%  '

$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").Value = '79.24'
$ws.Range("E39").Value = '  +10.66%  '

$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").Value = '0.585'
$ws.Range("E40").Value = '  +5.72%  '

$ws.Range("E41").Value = '  +1.20%  '

$ws.Range("D42").Value = '0.852'
$ws.Range("E42").Value = '  +2.83%  '

$ws.Range("E43").Value = '  +1.71%  '

$ws.Range("D44").Value = '0.0504'
$ws.Range("E44").Value = '  +1.26%  '

$ws.Range("D45").Value = '0.997'
$ws.Range("E45").Value = '  -0.30%  '

$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -4.29%  '

$ws.Range("D47").Value = '52.12'
$ws.Range("E47").Value = '  -5.71%  '

$ws.Range("D48").Value = '1.819.49'
$ws.Range("E48").Value = '  +2.81%  '

$ws.Range("D49").Value = '5.42'
$ws.Range("E49").Value = '  -0.53%  '

$ws.Range("D50").Value = '95.20'
$ws.Range("E50").Value = '  +6.25%  '

$ws.Range("E51").Value = '  +6.19%  '
